$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21
$ws.Range("B21").Value = "Sudden team member withdrawal "
$ws.Range("C21").Value = "Task delays and increased pressure on remaining members "
$ws.Range("D21").Value = "Cross-train members and document all work clearly "
$ws.Range("E21").Value = "High"
$ws.Range("E21").ClearFormats()
$ws.Range("E21").Font.Color = 192
$ws.Range("F21").Value = "Marina"
$ws.Range("G21").Value = "maintain updated tasks documentation for quick handover"

# Row 22
$ws.Range("B22").Value = "Environment differences between dev and deployment systems "
$ws.Range("C22").Value = "Unexpected bugs during deployment "
$ws.Range("D22").Value = "Set up a staging environment that mirrors production "
$ws.Range("E22").Value = "Medium"
$ws.Range("E22").ClearFormats()
$ws.Range("E22").Font.Color = 192
$ws.Range("F22").Value = "Marina"
$ws.Range("G22").Value = "Test deployment regularly before final release"

# Row 23
$ws.Range("B23").Value = "Lack of stakeholder feedback during development"
$ws.Range("C23").Value = "Misaligned product with user needs"
$ws.Range("D23").Value = "Schedule regular demos and feedback sessions with stakeholders"
$ws.Range("E23").Value = "Medium"
$ws.Range("E23").ClearFormats()
$ws.Range("E23").Font.Color = 192
$ws.Range("F23").Value = "Marina"
$ws.Range("G23").Value = "Use forms/surveys to gather structured feedback"

# Row 24
$ws.Range("B24").Value = "Incomplete testing coverage"
$ws.Range("C24").Value = "Undetected bugs and reliability issues"
$ws.Range("D24").Value = "Create and maintain test cases for all features, use test coverage tools"
$ws.Range("E24").Value = "High"
$ws.Range("E24").ClearFormats()
$ws.Range("E24").Font.Color = 192
$ws.Range("G24").Value = "Include test metrics in sprint review"
$ws.Range("F24").Value = "Marina "

# Widen column G and give custom widths to columns H and I
$ws.Columns("G").ColumnWidth = 48.666666666666664
$ws.Columns("H").ColumnWidth = 22.0
$ws.Columns("I").ColumnWidth = 32.0

# Update the active selection to B21 as in the edited file
$ws.Range("B21").Select()
